$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")
$ws.Select()
Write-Host "Sheet name: $($ws.Name)"
